$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.5
$ws.Range("H2").Value = 4.1
$ws.Range("I2").Value = 7
$ws.Range("J2").Value = 2.1
$ws.Range("L2").Value = 7
$ws.Range("M2").Value = 1.07
$ws.Range("N2").Value = 9
$ws.Range("U2").Value = 2.2
$ws.Range("V2").Value = 1.62
$ws.Range("Y2").Value = 8.5
$ws.Range("Z2").Value = 10
$ws.Range("AD2").Value = 8
$ws.Range("AG2").Value = 13
$ws.Range("AL2").Value = 51
$ws.Range("AN2").Value = 3.25
$ws.Range("AU2").Value = 9.5
$ws.Range("AV2").Value = 67
$ws.Range("AW2").Value = 8
$ws.Range("AZ2").Value = 151
$ws.Range("BB2").Value = 401
